# Update the destinations data (KS4 / KS5 sustained positive destination rate
# rows) to the latest available cohort (AY22/23 instead of AY21/22), and bump
# the "data period" label for those two rows from AY22/23 -> AY23/24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20: sustainedPositiveDestinationKS4Rate
$ws.Range("B20").Value = "AY23/24 data"
$ws.Range("C20").Value = "Destination measures show the percentage of students going to or remaining in an education, apprenticeship or employment destination in the academic year after completing Key Stage 4 studies (usually aged between 14 to 16). The cohort of learners used in the metrics here completed in AY22/23."

# Row 21: sustainedPositiveDestinationKS5Rate
$ws.Range("B21").Value = "AY23/24 data"
$ws.Range("C21").Value = "Destination measures show the percentage of students going to or remaining in an education, apprenticeship or employment destination in the academic year after completing Key Stage 5 studies (usually aged 18). The cohort of learners used in the metrics here completed in AY22/23."

# Scroll/select to reflect where the editor ended up after making the change.
$ws.Application.Goto($ws.Range("A22"), $true)
$ws.Range("B22").Select()
